$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40..129 down to 41..130
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with this week's new record. Columns that are
# identical across the whole table (A,B,C,E,F,G,H,I,N,Q,R) are copied from
# the row immediately below (the row that used to be row 40, now row 41),
# using Value2 to read (Value has a read quirk in this host for rvalue use).
$ws.Range("A40").Value = $ws.Range("A41").Value2
$ws.Range("B40").Value = $ws.Range("B41").Value2
$ws.Range("C40").Value = $ws.Range("C41").Value2
$ws.Range("D40").Value = 45014
$ws.Range("E40").Value = $ws.Range("E41").Value2
$ws.Range("F40").Value = $ws.Range("F41").Value2
$ws.Range("G40").Value = $ws.Range("G41").Value2
$ws.Range("H40").Value = $ws.Range("H41").Value2
$ws.Range("I40").Value = $ws.Range("I41").Value2
$ws.Range("J40").Value = 440
$ws.Range("K40").Value = 34000
$ws.Range("L40").Value = 35000
$ws.Range("M40").Value = 34500
$ws.Range("N40").Value = $ws.Range("N41").Value2
$ws.Range("O40").Value = "Provincia del Elquí"
$ws.Range("P40").Value = 1380
$ws.Range("Q40").Value = $ws.Range("Q41").Value2
$ws.Range("R40").Value = $ws.Range("R41").Value2

# Match the date cell format used by the rest of column D
$ws.Range("D40").NumberFormat = $ws.Range("D41").NumberFormat
